# Auto-generated update of market price columns (H:N) across multiple sheets
# Mirrors a scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 82
$ws.Range("I8").Value = 82
$ws.Range("K8").Value = 246
$ws.Range("M8").Value = -107

$ws.Range("H12").Value = 1377.7142
$ws.Range("I12").Value = 199
$ws.Range("J12").Value = 4324.5
$ws.Range("K12").Value = 199
$ws.Range("L12").Value = 4324.5
$ws.Range("M12").Value = -29
$ws.Range("N12").Value = -4664.5

$ws.Range("H15").Value = 599.89655
$ws.Range("I15").Value = 599.89655
$ws.Range("K15").Value = 1799.68965
$ws.Range("M15").Value = -1630.68965

$ws.Range("H18").Value = 2749
$ws.Range("I18").Value = 2749
$ws.Range("K18").Value = 2749
$ws.Range("M18").Value = -2465

$ws.Range("H33").Value = 256.75
$ws.Range("I33").Value = 231.15384
$ws.Range("J33").Value = 367.66666
$ws.Range("K33").Value = 231.15384
$ws.Range("L33").Value = 367.66666
$ws.Range("M33").Value = -2.153840000000002
$ws.Range("N33").Value = -825.66666

$ws.Range("H49").Value = 1000
$ws.Range("J49").Value = 1000
$ws.Range("L49").Value = 3000
$ws.Range("N49").Value = -3272

$ws.Range("H88").Value = 1967.4762
$ws.Range("J88").Value = 2156.611
$ws.Range("L88").Value = 2156.611
$ws.Range("N88").Value = -2968.611

$ws.Range("H91").Value = 1967.4762
$ws.Range("J91").Value = 2156.611
$ws.Range("L91").Value = 2156.611
$ws.Range("N91").Value = -4964.611

$ws.Range("H125").Value = 50002604
$ws.Range("I125").Value = 83334230
$ws.Range("J125").Value = 5165.5
$ws.Range("K125").Value = 750008070
$ws.Range("L125").Value = 46489.5
$ws.Range("M125").Value = -750005610
$ws.Range("N125").Value = -51409.5

$ws.Range("H127").Value = 800.6667
$ws.Range("I127").Value = 701
$ws.Range("K127").Value = 2103
$ws.Range("M127").Value = 2857

$ws.Range("H131").Value = 1057.8
$ws.Range("I131").Value = 822.25
$ws.Range("K131").Value = 2466.75
$ws.Range("M131").Value = 2573.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 409.36365
$ws.Range("J5").Value = 393.66666
$ws.Range("L5").Value = 393.66666
$ws.Range("N5").Value = -617.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 409.36365
$ws.Range("J4").Value = 393.66666
$ws.Range("L4").Value = 393.66666
$ws.Range("N4").Value = -623.66666

$ws.Range("H37").Value = 1761.2
$ws.Range("I37").Value = 951.5
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 951.5
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -814.5
$ws.Range("N37").Value = -5274

$ws.Range("H58").Value = 51425
$ws.Range("J58").Value = 51425
$ws.Range("L58").Value = 51425
$ws.Range("N58").Value = -52013

$ws.Range("H107").Value = 1785
$ws.Range("I107").Value = 1785
$ws.Range("K107").Value = 1785
$ws.Range("M107").Value = 135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1728.1666
$ws.Range("J7").Value = 2093.6
$ws.Range("L7").Value = 2093.6
$ws.Range("N7").Value = -2319.6

$ws.Range("H22").Value = 1455.4445
$ws.Range("J22").Value = 1624.875
$ws.Range("L22").Value = 1624.875
$ws.Range("N22").Value = -2324.875

$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -47290

$ws.Range("H60").Value = 25998.666
$ws.Range("J60").Value = 25998.666
$ws.Range("L60").Value = 25998.666
$ws.Range("N60").Value = -27020.666

$ws.Range("H68").Value = 47648.832
$ws.Range("J68").Value = 47648.832
$ws.Range("L68").Value = 47648.832
$ws.Range("N68").Value = -49146.832

$ws.Range("H71").Value = 47648.832
$ws.Range("J71").Value = 47648.832
$ws.Range("L71").Value = 142946.496
$ws.Range("N71").Value = -150434.496

$ws.Range("H99").Value = 6117.4
$ws.Range("I99").Value = 4886.1816
$ws.Range("K99").Value = 4886.1816
$ws.Range("M99").Value = -3388.1816

$ws.Range("H126").Value = 6117.4
$ws.Range("I126").Value = 4886.1816
$ws.Range("K126").Value = 14658.5448
$ws.Range("M126").Value = -12188.5448

$ws.Range("H141").Value = 387499.75
$ws.Range("J141").Value = 416666.34
$ws.Range("L141").Value = 416666.34
$ws.Range("N141").Value = -427026.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 138.8
$ws.Range("I33").Value = 198.5
$ws.Range("J33").Value = 99
$ws.Range("K33").Value = 1191
$ws.Range("L33").Value = 594
$ws.Range("M33").Value = -908
$ws.Range("N33").Value = -1160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1654.6666
$ws.Range("I102").Value = 1584.6
$ws.Range("K102").Value = 1584.6
$ws.Range("M102").Value = 37.40000000000009

$ws.Range("H137").Value = 70189.5
$ws.Range("J137").Value = 70189.5
$ws.Range("L137").Value = 70189.5
$ws.Range("N137").Value = -80389.5

$ws.Range("H138").Value = 65333
$ws.Range("J138").Value = 65333
$ws.Range("L138").Value = 65333
$ws.Range("N138").Value = -75613

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 994.25
$ws.Range("I16").Value = 992.3333
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 992.3333
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -822.3333
$ws.Range("N16").Value = -1340

$ws.Range("H22").Value = 2019
$ws.Range("I22").Value = 1865
$ws.Range("J22").Value = 2250
$ws.Range("K22").Value = 1865
$ws.Range("L22").Value = 2250
$ws.Range("M22").Value = -1570
$ws.Range("N22").Value = -2840

$ws.Range("H27").Value = 2019
$ws.Range("I27").Value = 1865
$ws.Range("J27").Value = 2250
$ws.Range("K27").Value = 1865
$ws.Range("L27").Value = 2250
$ws.Range("M27").Value = -1758
$ws.Range("N27").Value = -2464

$ws.Range("H46").Value = 2167.818
$ws.Range("I46").Value = 1654.6666
$ws.Range("J46").Value = 2523.077
$ws.Range("K46").Value = 1654.6666
$ws.Range("L46").Value = 2523.077
$ws.Range("M46").Value = -1466.6666
$ws.Range("N46").Value = -2899.077

$ws.Range("H55").Value = 939.85
$ws.Range("J55").Value = 1035.1177
$ws.Range("L55").Value = 1035.1177
$ws.Range("N55").Value = -1381.1177

$ws.Range("H100").Value = 3998.4285
$ws.Range("J100").Value = 3998.4285
$ws.Range("L100").Value = 3998.4285
$ws.Range("N100").Value = -5080.4285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1388.7273
$ws.Range("I126").Value = 1388.7273
$ws.Range("K126").Value = 4166.1819
$ws.Range("M126").Value = -1696.1819

$ws.Range("H132").Value = 127083.125
$ws.Range("I132").Value = 127083.125
$ws.Range("K132").Value = 381249.375
$ws.Range("M132").Value = -378719.375

$ws.Range("H136").Value = 1308
$ws.Range("I136").Value = 961.9394
$ws.Range("J136").Value = 3592
$ws.Range("K136").Value = 2885.8182
$ws.Range("L136").Value = 10776
$ws.Range("M136").Value = -335.8181999999997
$ws.Range("N136").Value = -15876
